# Refresh the crypto price/volume snapshot (GitHub Actions scheduled update).
# Each value is written with a leading apostrophe so Excel stores it as text
# (matching the original inlineStr cells) instead of auto-converting it to a
# number/date.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'70.778.22"
$ws.Range("E2").Value = "'  +2.46%  "
$ws.Range("D3").Value = "'3.562.66"
$ws.Range("E3").Value = "'  +1.46%  "
$ws.Range("E4").Value = "'  -0.04%  "
$ws.Range("D5").Value = "'580.85"
$ws.Range("E5").Value = "'  +1.44%  "
$ws.Range("D6").Value = "'187.27"
$ws.Range("E6").Value = "'  +1.51%  "
$ws.Range("E7").Value = "'  +2.17%  "
$ws.Range("D8").Value = "'3.553.67"
$ws.Range("E8").Value = "'  +1.46%  "
$ws.Range("E9").Value = "'  -0.06%  "
$ws.Range("D10").Value = "'0.219"
$ws.Range("E10").Value = "'  +18.85%  "
$ws.Range("D11").Value = "'0.651"
$ws.Range("E11").Value = "'  -0.09%  "
$ws.Range("D12").Value = "'54.49"
$ws.Range("E12").Value = "'  +0.43%  "
$ws.Range("E13").Value = "'  +5.89%  "
$ws.Range("D14").Value = "'9.52"
$ws.Range("E14").Value = "'  +0.83%  "
$ws.Range("D15").Value = "'4.133.08"
$ws.Range("D16").Value = "'70.835.78"
$ws.Range("E16").Value = "'  +2.59%  "
$ws.Range("D17").Value = "'12.84"
$ws.Range("E17").Value = "'  +4.29%  "
$ws.Range("D18").Value = "'19.19"
$ws.Range("E18").Value = "'  -1.06%  "
$ws.Range("D19").Value = "'3.566.87"
$ws.Range("E19").Value = "'  +1.83%  "
$ws.Range("D20").Value = "'572.90"
$ws.Range("E20").Value = "'  +5.56%  "
$ws.Range("E21").Value = "'  +0.69%  "
$ws.Range("E22").Value = "'  -0.89%  "
$ws.Range("D23").Value = "'17.65"
$ws.Range("E23").Value = "'  -3.71%  "
$ws.Range("E24").Value = "'  +3.68%  "
$ws.Range("E25").Value = "'  -2.04%  "
$ws.Range("D26").Value = "'94.34"
$ws.Range("E26").Value = "'  +0.19%  "
$ws.Range("D27").Value = "'11.20"
$ws.Range("E27").Value = "'  +1.64%  "
$ws.Range("E28").Value = "'  +1.56%  "
$ws.Range("E29").Value = "'  +2.76%  "
$ws.Range("D30").Value = "'32.76"
$ws.Range("E30").Value = "'  +3.33%  "
$ws.Range("D31").Value = "'7.21"
$ws.Range("E31").Value = "'  -0.62%  "
$ws.Range("D32").Value = "'12.32"
$ws.Range("E32").Value = "'  -2.80%  "
$ws.Range("E34").Value = "'  +23.23%  "
$ws.Range("D35").Value = "'63.36"
$ws.Range("E35").Value = "'  -1.88%  "
$ws.Range("D36").Value = "'3.29"
$ws.Range("E36").Value = "'  +6.89%  "
$ws.Range("D37").Value = "'536.27"
$ws.Range("E37").Value = "'  -3.36%  "
$ws.Range("D38").Value = "'0.411"
$ws.Range("E38").Value = "'  +2.70%  "
$ws.Range("D39").Value = "'0.0₃0811"
$ws.Range("E39").Value = "'  +6.29%  "
$ws.Range("D40").Value = "'38.10"
$ws.Range("E40").Value = "'  +0.22%  "
$ws.Range("E41").Value = "'  +0.00%  "
$ws.Range("D42").Value = "'3.632.79"
$ws.Range("E42").Value = "'  +10.70%  "
$ws.Range("E43").Value = "'  +5.30%  "
$ws.Range("E44").Value = "'  +2.78%  "
$ws.Range("D45").Value = "'0.0471"
$ws.Range("E45").Value = "'  +5.86%  "
$ws.Range("B46").Value = "'ApeXProtocol"
$ws.Range("C46").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'3.46"
$ws.Range("E46").Value = "'  -0.22%  "
$ws.Range("B47").Value = "'ThetaToken"
$ws.Range("C47").Value = "'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").Value = "'2.94"
$ws.Range("E47").Value = "'  -1.87%  "
$ws.Range("E48").Value = "'  +5.17%  "
$ws.Range("E49").Value = "'  +3.03%  "
$ws.Range("D50").Value = "'1.00"
$ws.Range("E50").Value = "'  +0.03%  "
$ws.Range("E51").Value = "'  +5.49%  "
